# "new server version 2.0" - update strategy parameters (Take Profit 2 / Stop loss
# values) in the trading-pairs worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ETHUSDT
$ws.Range("E2").Value = 0.1
$ws.Range("G2").Value = 0.25
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.05

# Row 3 - LINKUSDT
$ws.Range("E3").Value = 0.1
$ws.Range("G3").Value = 0.5
$ws.Range("I3").Value = 0.5
$ws.Range("J3").Value = 0.05

# Row 12 - OMGUSDT
$ws.Range("E12").Value = 0.2
$ws.Range("G12").Value = 0.8
$ws.Range("I12").Value = 0.9
$ws.Range("J12").Value = 0.05

# Row 13 - SLPUSDT
$ws.Range("G13").Value = 0.4

# Row 14 - RUNEUSDT
$ws.Range("G14").Value = 0.4

# Row 15 - BTCUSDT
$ws.Range("G15").Value = 0.4

# Row 16 - SRMUSDT
$ws.Range("G16").Value = 0.4

# Row 17 - FILUSDT
$ws.Range("G17").Value = 0.4

# Leave the selection where the author finished editing
[void]$ws.Range("J4").Select()
